# Brewing Potions.xlsx — apply commit changes via Excel COM interop
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Potions" (sheet1)
# ---------------------------------------------------------------------------
$potions = $wb.Worksheets.Item("Potions")
$potions.Activate()

# Row 42: bold the row (new font + reuses existing fills/border -> new cellXfs 8/9)
$potions.Range("A42:G42").Font.Bold = $true
$potions.Range("I42").Font.Bold = $true
$potions.Range("K42").Font.Bold = $true
$potions.Range("H42").Font.Bold = $true
$potions.Range("J42").Font.Bold = $true

# Row 50: new potion entry ("Tongues")
$potions.Range("A50").Value = "Tongues"
$potions.Range("B50").Value = 3
$potions.Range("C50").Value = 3
$potions.Range("D50").Value = 10
$potions.Range("E50").Value = "Min"
$potions.Range("F50").Formula = "=(B50*C50*50)"
$potions.Range("G50").Formula = "=F50-(F50*0.05)"
$potions.Range("H50").Formula = "=G50/2"
$potions.Range("I50").Value = 8
$potions.Range("J50").Formula = "=I50/2"
$potions.Range("K50").Formula = "=5+C50"

# Rows 52-58: fill in column A with "Cure Serious Wound"
$potions.Range("A52").Value = "Cure Serious Wound"
$potions.Range("A53").Value = "Cure Serious Wound"
$potions.Range("A54").Value = "Cure Serious Wound"
$potions.Range("A55").Value = "Cure Serious Wound"
$potions.Range("A56").Value = "Cure Serious Wound"
$potions.Range("A57").Value = "Cure Serious Wound"
$potions.Range("A58").Value = "Cure Serious Wound"

# sheet view: scroll + selection
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$potions.Range("H32").Select()

# ---------------------------------------------------------------------------
# Sheet "Poisons" (sheet2)
# ---------------------------------------------------------------------------
$poisons = $wb.Worksheets.Item("Poisons")
$poisons.Activate()
$excel.ActiveWindow.ScrollRow = 82
$excel.ActiveWindow.ScrollColumn = 1
$poisons.Range("A34").Select()

# restore original active sheet/tab
$potions.Activate()
